$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells receive plain numeric-looking text (e.g. "560.83") that Excel
# would otherwise auto-convert to a floating point number. Mark them as Text
# first so the literal string from the source data is preserved exactly.
$textCellRefs = @('D5', 'D6', 'D7', 'D9', 'D11', 'D12', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D44', 'D49')
foreach ($ref in $textCellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.366.77'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.100.59'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '560.83'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = '144.47'
$ws.Range('E6').Value = '  +2.49%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.100.72'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').Value = '6.17'
$ws.Range('E11').Value = '  -4.91%  '
$ws.Range('D12').Value = '0.473'
$ws.Range('E12').Value = '  +3.57%  '
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').Value = '35.23'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '3.600.36'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '64.382.29'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').Value = '3.097.90'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '0.111'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').Value = '484.39'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '14.02'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '0.679'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '7.58'
$ws.Range('E23').Value = '  +4.08%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '14.33'
$ws.Range('E24').Value = '  +12.39%  '
$ws.Range('D25').Value = '81.43'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '2.81'
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('D28').Value = '8.04'
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '26.46'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').Value = '5.64'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').Value = '6.25'
$ws.Range('E35').Value = '  +4.13%  '
$ws.Range('D36').Value = '55.63'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '2.99'
$ws.Range('E37').Value = '  +16.02%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0410'
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').Value = '447.68'
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('D40').Value = '0.0818'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').Value = '2.969.07'
$ws.Range('E41').Value = '  -2.81%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -4.97%  '
$ws.Range('D44').Value = '28.25'
$ws.Range('E44').Value = '  +1.00%  '
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('E47').Value = '  +4.57%  '
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').Value = '118.53'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('E51').Value = '  +0.29%  '
